# Usman - users added to run files
#
# Sheet3 previously only flagged a single test run with "Y" in cell A2.
# Extend that flag down through A8 (rows 3-8) so the additional users
# added to the run files are included.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

$ws.Range("A3:A8").Value = "Y"
